$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 311, pushing the existing rows
# (previously 311-428) down to 313-430.
$ws.Rows("311:312").Insert()

# New row 311: Start Ruby / Primera, date 2022-11-11 (serial 44876)
$ws.Range("A311").Value = 4
$ws.Range("B311").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C311").Value = "Los Lagos"
$ws.Range("D311").Value = 44876
$ws.Range("E311").Value = 10
$ws.Range("F311").Value = "Fruta"
$ws.Range("G311").Value = 100102
$ws.Range("H311").Value = "Cítricos"
$ws.Range("I311").Value = 100102006
$ws.Range("J311").Value = "Pomelo"
$ws.Range("K311").Value = "Start Ruby"
$ws.Range("L311").Value = "Primera"
$ws.Range("M311").Value = 200
$ws.Range("N311").Value = 14000
$ws.Range("O311").Value = 15000
$ws.Range("P311").Value = 14500
$ws.Range("Q311").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R311").Value = "Región de O'Higgins"
$ws.Range("S311").Value = 1036
$ws.Range("T311").Value = 14

# New row 312: Start Ruby / Segunda, same date 2022-11-11 (serial 44876)
$ws.Range("A312").Value = 4
$ws.Range("B312").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C312").Value = "Los Lagos"
$ws.Range("D312").Value = 44876
$ws.Range("E312").Value = 10
$ws.Range("F312").Value = "Fruta"
$ws.Range("G312").Value = 100102
$ws.Range("H312").Value = "Cítricos"
$ws.Range("I312").Value = 100102006
$ws.Range("J312").Value = "Pomelo"
$ws.Range("K312").Value = "Start Ruby"
$ws.Range("L312").Value = "Segunda"
$ws.Range("M312").Value = 100
$ws.Range("N312").Value = 12000
$ws.Range("O312").Value = 12000
$ws.Range("P312").Value = 12000
$ws.Range("Q312").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R312").Value = "Región de O'Higgins"
$ws.Range("S312").Value = 857
$ws.Range("T312").Value = 14
